$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $c = $ws.Range($cellRef)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue "D2" "41.869.60"
Set-TextValue "E2" "  -0.19%  "
Set-TextValue "D3" "2.255.66"
Set-TextValue "E3" "  -0.90%  "
Set-TextValue "E4" "  +0.12%  "
Set-TextValue "D5" "304.45"
Set-TextValue "E5" "  -0.35%  "
Set-TextValue "D6" "95.04"
Set-TextValue "E6" "  +1.96%  "
Set-TextValue "D7" "0.524"
Set-TextValue "E7" "  -1.19%  "
Set-TextValue "E8" "  +0.14%  "
Set-TextValue "D9" "0.486"
Set-TextValue "E9" "  -0.38%  "
Set-TextValue "D10" "34.59"
Set-TextValue "E10" "  +5.59%  "
Set-TextValue "D11" "0.0785"
Set-TextValue "E11" "  -1.67%  "
Set-TextValue "E12" "  -0.39%  "
Set-TextValue "D13" "6.60"
Set-TextValue "E13" "  -1.54%  "
Set-TextValue "D14" "2.614.11"
Set-TextValue "E14" "  -0.43%  "
Set-TextValue "D15" "14.27"
Set-TextValue "E15" "  -0.80%  "
Set-TextValue "D16" "2.263.41"
Set-TextValue "E16" "  -0.52%  "
Set-TextValue "D17" "0.787"
Set-TextValue "E17" "  +0.43%  "
Set-TextValue "D18" "41.778.03"
Set-TextValue "E18" "  -0.20%  "
Set-TextValue "D19" "12.25"
Set-TextValue "E19" "  -4.68%  "
Set-TextValue "D20" "0.0₃0897"
Set-TextValue "E20" "  -2.17%  "
Set-TextValue "D21" "5.92"
Set-TextValue "E21" "  -1.05%  "
Set-TextValue "D22" "67.78"
Set-TextValue "E22" "  -0.35%  "
Set-TextValue "D23" "236.21"
Set-TextValue "E23" "  -3.23%  "
Set-TextValue "D24" "2.55"
Set-TextValue "E24" "  -2.42%  "
Set-TextValue "E25" "  -0.11%  "
Set-TextValue "D26" "1.91"
Set-TextValue "E26" "  -1.48%  "
Set-TextValue "D27" "23.54"
Set-TextValue "E27" "  -2.14%  "
Set-TextValue "D28" "36.33"
Set-TextValue "E28" "  +4.20%  "
Set-TextValue "D29" "2.11"
Set-TextValue "E29" "  +1.41%  "
Set-TextValue "D30" "9.43"
Set-TextValue "E30" "  -2.79%  "
Set-TextValue "D31" "159.69"
Set-TextValue "E31" "  +0.34%  "
Set-TextValue "B32" "FirstDigitalUSD"
Set-TextValue "C32" "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
Set-TextValue "D32" "1.00"
Set-TextValue "E32" "  +0.04%  "
Set-TextValue "B33" "Filecoin"
Set-TextValue "C33" "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
Set-TextValue "D33" "5.17"
Set-TextValue "E33" "  -3.99%  "
Set-TextValue "D34" "3.12"
Set-TextValue "E34" "  +2.77%  "
Set-TextValue "D35" "0.0732"
Set-TextValue "E35" "  -1.63%  "
Set-TextValue "D36" "17.01"
Set-TextValue "E36" "  -2.17%  "
Set-TextValue "D37" "2.38"
Set-TextValue "E37" "  +0.53%  "
Set-TextValue "E38" "  -1.20%  "
Set-TextValue "D39" "1.81"
Set-TextValue "E39" "  +0.69%  "
Set-TextValue "D40" "0.113"
Set-TextValue "E40" "  -2.76%  "
Set-TextValue "D41" "3.98"
Set-TextValue "E41" "  +1.05%  "
Set-TextValue "D42" "2.37"
Set-TextValue "E42" "  +5.36%  "
Set-TextValue "D43" "1.971.53"
Set-TextValue "E43" "  -1.85%  "
Set-TextValue "D44" "0.0282"
Set-TextValue "E44" "  -0.31%  "
Set-TextValue "D45" "18.48"
Set-TextValue "E45" "  -6.68%  "
Set-TextValue "D46" "2.91"
Set-TextValue "E46" "  -0.41%  "
Set-TextValue "D47" "9.82"
Set-TextValue "E47" "  -4.47%  "
Set-TextValue "D48" "53.04"
Set-TextValue "E48" "  -1.01%  "
Set-TextValue "D49" "72.38"
Set-TextValue "E49" "  -0.54%  "
Set-TextValue "E50" "  -0.70%  "
Set-TextValue "D51" "90.22"
Set-TextValue "E51" "  -1.94%  "
